$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 5 -> 4, Wrong penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 80 -> 64, Wrong penalty total -2 -> -4, Max text "80 / 140" -> "60 / 112"
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "60 / 112"
